# Update "想去人数" (want-to-go count) figures on the "展览" and "全部类型" sheets
# to reflect the refreshed data pulled at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F4").Value = 634
$wsExhibition.Range("F7").Value = 2771
$wsExhibition.Range("F9").Value = 7770
$wsExhibition.Range("F12").Value = 41
$wsExhibition.Range("F13").Value = 331

# --- Sheet "全部类型" ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 634
$wsAll.Range("F9").Value = 2771
$wsAll.Range("F11").Value = 7770
$wsAll.Range("F14").Value = 41
$wsAll.Range("F17").Value = 331
